# Generate Report for Handback
# Update "Latest Handback DateTime" (column L) for the 8fc3a59f-... row (row 3)
# on both the zh-cn and de-de localization status sheets, reflecting a newly
# generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("L3").Value = "2017-01-03 06:15:28"
$wsZhCn.Range("L3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("L3").Value = "2017-01-03 06:15:39"
$wsDeDe.Range("L3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
